$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 645.7027
$ws.Range("J17").Value = 535.10767
$ws.Range("L17").Value = 1605.32301
$ws.Range("N17").Value = -1941.32301

$ws.Range("H21").Value = 33892
$ws.Range("I21").Value = 35005.285
$ws.Range("J21").Value = 29995.5
$ws.Range("K21").Value = 35005.285
$ws.Range("L21").Value = 29995.5
$ws.Range("M21").Value = -34537.285
$ws.Range("N21").Value = -30931.5

$ws.Range("H23").Value = 33892
$ws.Range("I23").Value = 35005.285
$ws.Range("J23").Value = 29995.5
$ws.Range("K23").Value = 35005.285
$ws.Range("L23").Value = 29995.5
$ws.Range("M23").Value = -34771.285
$ws.Range("N23").Value = -30463.5

$ws.Range("H28").Value = 617.6896400000001
$ws.Range("J28").Value = 595.75
$ws.Range("L28").Value = 595.75
$ws.Range("N28").Value = -1565.75

$ws.Range("H86").Value = 1614.3077
$ws.Range("J86").Value = 1813.8572
$ws.Range("L86").Value = 1813.8572
$ws.Range("N86").Value = -4059.8572

$ws.Range("H89").Value = 1614.3077
$ws.Range("J89").Value = 1813.8572
$ws.Range("L89").Value = 9069.286
$ws.Range("N89").Value = -20301.286

$ws.Range("H92").Value = 1659.579
$ws.Range("I92").Value = 1753.6923
$ws.Range("J92").Value = 1455.6666
$ws.Range("K92").Value = 1753.6923
$ws.Range("L92").Value = 1455.6666
$ws.Range("M92").Value = -505.6922999999999
$ws.Range("N92").Value = -3951.6666

$ws.Range("H129").Value = 917.798
$ws.Range("I129").Value = 471
$ws.Range("J129").Value = 951.79346
$ws.Range("K129").Value = 1413
$ws.Range("L129").Value = 2855.38038
$ws.Range("M129").Value = 3587
$ws.Range("N129").Value = -12855.38038

$ws.Range("H132").Value = 175326.97
$ws.Range("I132").Value = 2773.18
$ws.Range("K132").Value = 8319.539999999999
$ws.Range("M132").Value = -5789.539999999999

$ws.Range("H137").Value = 2373.6428
$ws.Range("I137").Value = 1311.1482
$ws.Range("J137").Value = 4286.1333
$ws.Range("K137").Value = 3933.4446
$ws.Range("L137").Value = 12858.3999
$ws.Range("M137").Value = -1383.4446
$ws.Range("N137").Value = -17958.3999

$ws.Range("H138").Value = 4979.87
$ws.Range("J138").Value = 7031.8506
$ws.Range("L138").Value = 21095.5518
$ws.Range("N138").Value = -31375.5518

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5786.387
$ws.Range("I32").Value = 4881.9614
$ws.Range("J32").Value = 10489.4
$ws.Range("K32").Value = 4881.9614
$ws.Range("L32").Value = 10489.4
$ws.Range("M32").Value = -4594.9614
$ws.Range("N32").Value = -11063.4

$ws.Range("H132").Value = 2454.7307
$ws.Range("I132").Value = 1101.5294
$ws.Range("J132").Value = 5010.778
$ws.Range("K132").Value = 3304.5882
$ws.Range("L132").Value = 15032.334
$ws.Range("M132").Value = -774.5881999999997
$ws.Range("N132").Value = -20092.334

$ws.Range("H139").Value = 43080
$ws.Range("J139").Value = 43080
$ws.Range("L139").Value = 43080
$ws.Range("N139").Value = -53360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H43").Value = 79800
$ws.Range("J43").Value = 79800
$ws.Range("L43").Value = 79800
$ws.Range("N43").Value = -80162

$ws.Range("H134").Value = 1593.1702
$ws.Range("I134").Value = 1118.2222
$ws.Range("J134").Value = 3147.5454
$ws.Range("K134").Value = 3354.6666
$ws.Range("L134").Value = 9442.636200000001
$ws.Range("M134").Value = -819.6665999999996
$ws.Range("N134").Value = -14512.6362

$ws.Range("H138").Value = 41305.715
$ws.Range("J138").Value = 41305.715
$ws.Range("L138").Value = 41305.715
$ws.Range("N138").Value = -51585.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14287954
$ws.Range("I31").Value = 1006.381
$ws.Range("J31").Value = 35718376
$ws.Range("K31").Value = 1006.381
$ws.Range("L31").Value = 35718376
$ws.Range("M31").Value = -711.381
$ws.Range("N31").Value = -35718966

$ws.Range("H34").Value = 14287954
$ws.Range("I34").Value = 1006.381
$ws.Range("J34").Value = 35718376
$ws.Range("K34").Value = 1006.381
$ws.Range("L34").Value = 35718376
$ws.Range("M34").Value = -804.381
$ws.Range("N34").Value = -35718780

$ws.Range("H58").Value = 1569.2954
$ws.Range("I58").Value = 1478.5363
$ws.Range("J58").Value = 1898.8948
$ws.Range("K58").Value = 1478.5363
$ws.Range("L58").Value = 1898.8948
$ws.Range("M58").Value = -1275.5363
$ws.Range("N58").Value = -2304.8948

$ws.Range("H62").Value = 10000
$ws.Range("I62").Value = 10000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 10000
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -9376

$ws.Range("H65").Value = 10000
$ws.Range("I65").Value = 10000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 50000
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -46880

$ws.Range("H134").Value = 7595.9443
$ws.Range("I134").Value = 8517.923000000001
$ws.Range("J134").Value = 5198.8
$ws.Range("K134").Value = 25553.769
$ws.Range("L134").Value = 15596.4
$ws.Range("M134").Value = -23018.769
$ws.Range("N134").Value = -20666.4

$ws.Range("H136").Value = 1569.2954
$ws.Range("I136").Value = 1478.5363
$ws.Range("J136").Value = 1898.8948
$ws.Range("K136").Value = 4435.6089
$ws.Range("L136").Value = 5696.6844
$ws.Range("M136").Value = -1885.6089
$ws.Range("N136").Value = -10796.6844

$ws.Range("H138").Value = 42733
$ws.Range("J138").Value = 42733
$ws.Range("L138").Value = 42733
$ws.Range("N138").Value = -53013

$ws.Range("H140").Value = 77527.336
$ws.Range("J140").Value = 77527.336
$ws.Range("L140").Value = 77527.336
$ws.Range("N140").Value = -87887.336

$ws.Range("H141").Value = 33037.5
$ws.Range("J141").Value = 33037.5
$ws.Range("L141").Value = 33037.5
$ws.Range("N141").Value = -43397.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 4000
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").Value = 0

$ws.Range("H140").Value = 38985.555
$ws.Range("J140").Value = 38985.555
$ws.Range("L140").Value = 38985.555
$ws.Range("N140").Value = -49345.555

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 537.5
$ws.Range("I9").Value = 150
$ws.Range("J9").Value = 666.6667
$ws.Range("K9").Value = 150
$ws.Range("L9").Value = 666.6667
$ws.Range("M9").Value = 74
$ws.Range("N9").Value = -1114.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 19660
$ws.Range("J24").Value = 27990
$ws.Range("L24").Value = 27990
$ws.Range("N24").Value = -28450

$ws.Range("H132").Value = 7409035.5
$ws.Range("I132").Value = 907.63635
$ws.Range("K132").Value = 2722.90905
$ws.Range("M132").Value = -192.9090500000002
